# Refresh the "kfranklin3" web query tracklist into the new "jayelectronica1"
# tracklist (Jay Electronica - A Written Testimony), per commit
# "11092220 added a review je".
#
# The workbook keeps two query-result sheets (Sheet1 and Sheet3, both fed by
# the same web query under two connections) plus Sheet2, a fully
# formula-driven "review" sheet that re-formats whatever is in Sheet1 - so
# only Sheet1/Sheet3 need new literal data; Sheet2 recalculates on its own.

$wb = $excel.ActiveWorkbook

# New tracklist: title / composers / performer / duration (as h:mm fraction
# of a day, matching what a refreshed web query would have written).
$tracks = @(
  @("The Overwhelming Event", "Eric Demarsan / Larrance Dopson / Louis Farrakhan / Khirye Tyler", "Jay Electronica", 0.094444444444444442),
  @("Ghost of Soulja Slim", "Elpadaro F. Electronica Allah / Shawn Carter / Larrance Dopson / Louis Farrakhan / Craig Lawson / Chris Payton / James Tapp / Khirye Tyler / John Williams", "Jay Electronica", 0.18472222222222223),
  @("The Blinding", "Elpadaro F. Electronica Allah / Jahron Brathwaite / Shawn Carter / Kaseem Dean / Chauncey Hollis / Abraham Orellana / Jacques Webster", "Jay Electronica", 0.11666666666666665),
  @("The Neverending Story", "Elpadaro F. Electronica Allah / Shawn Carter / Mirtha Defilpo / Daniel Alan Maman / Terius Nash / Felix Nebbia / Litto Nebbia", "Jay Electronica", 0.18055555555555555),
  @("Shiny Suit Theory", "Elpadaro F. Electronica Allah / Shawn Carter / Barbara Mason / Terius Nash", "Jay Electronica", 0.16874999999999998),
  @("Universal Soldier", "Elpadaro F. Electronica Allah / James Blake / Jennifer Vashti Bunyan / Shawn Carter / Allen Toussaint / Jacques Webster", "Jay Electronica", 0.17916666666666667),
  @("Flux Capacitor", "Elpadaro F. Electronica Allah / Badriia Bourelly / Jerry Butler / Shawn Carter / James Fauntleroy / Robyn Fenty / Kenny Gamble / Leon Huff / Elton Newman / Lawrence Parker / D`u2019Artanian Stovall / Dion Wilson", "Jay Electronica", 0.14305555555555557),
  @("Fruits of the Spirit", "Elpadaro F. Electronica Allah / Rodney G. Massey / Dion Wilson", "Jay Electronica", 0.065277777777777782),
  @("Ezekiel`u2019s Wheel", "Elpadaro F. Electronica Allah / Shawn Carter / Brian Eno / James Fauntleroy / Robert Fripp", "Jay Electronica", 0.28263888888888888),
  @("A.P.I.D.T.A.", "Elpadaro F. Electronica Allah / Shawn Carter / Laura Manders / Mark Speer", "Jay Electronica", 0.22777777777777777)
)

foreach ($sheetName in @("Sheet1", "Sheet3")) {
  $ws = $wb.Worksheets.Item($sheetName)

  # Row 2..11: write the new 10-track result set (was 13 tracks before).
  for ($i = 0; $i -lt $tracks.Count; $i++) {
    $row = $i + 2
    $track = $tracks[$i]
    $ws.Cells.Item($row, 1).Value2 = $i + 1
    $ws.Cells.Item($row, 2).Value2 = $track[0]
    $ws.Cells.Item($row, 3).Value2 = $track[1]
    $ws.Cells.Item($row, 4).Value2 = $track[2]
    $ws.Cells.Item($row, 5).Value2 = $track[3]
  }

  # Rows 12..14 held tracks 11..13 of the old 13-track result; the new query
  # only returned 10 rows, so the web-query refresh clears them out.
  $ws.Range("A12:E14").ClearContents()

  # A refreshed query also re-applies the result column widths.
  $ws.Columns.Item(2).ColumnWidth = 23.21875
  $ws.Columns.Item(3).ColumnWidth = 80.88671875
  $ws.Columns.Item(4).ColumnWidth = 14.109375
  $ws.Columns.Item(5).ColumnWidth = 5.6640625
}

# The query's defined name (used as the refresh range) gets renamed from
# "kfranklin3" to "jayelectronica1" and shrinks to the new 10-row extent.
$n1 = $wb.Names.Item("Sheet1!kfranklin3")
$n1.Name = "jayelectronica1"
$wb.Names.Item("Sheet1!jayelectronica1").RefersTo = '=Sheet1!$A$1:$E$11'

$n3 = $wb.Names.Item("Sheet3!kfranklin3")
$n3.Name = "jayelectronica1"
$wb.Names.Item("Sheet3!jayelectronica1").RefersTo = '=Sheet3!$A$1:$E$11'

# Sheet2 (the formatted review sheet) only has 10 data rows now, so the
# lingering selection shrinks from row 17 to row 14 to match.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("A3:K14").Select()
